$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 313
$ws.Range("I2").Value = 250
$ws.Range("K2").Value = 250
$ws.Range("M2").Value = -137
# Row 12
$ws.Range("H12").Value = 844.9
$ws.Range("I12").Value = 487.25
$ws.Range("J12").Value = 1083.3334
$ws.Range("K12").Value = 487.25
$ws.Range("L12").Value = 1083.3334
$ws.Range("M12").Value = -317.25
$ws.Range("N12").Value = -1423.3334
# Row 62
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
# Row 65
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
# Row 80
$ws.Range("H80").Value = 3253.6
$ws.Range("I80").Value = 3522.8333
$ws.Range("K80").Value = 10568.4999
$ws.Range("M80").Value = -9570.499899999999
# Row 83
$ws.Range("H83").Value = 3253.6
$ws.Range("I83").Value = 3522.8333
$ws.Range("K83").Value = 31705.4997
$ws.Range("M83").Value = -26713.4997
# Row 97
$ws.Range("H97").Value = 3900
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()
# Row 98
$ws.Range("H98").Value = 923.7143
$ws.Range("I98").Value = 818.5294
$ws.Range("K98").Value = 818.5294
$ws.Range("M98").Value = 679.4706
# Row 122
$ws.Range("H122").Value = 923.7143
$ws.Range("I122").Value = 818.5294
$ws.Range("K122").Value = 2455.5882
$ws.Range("M122").Value = -5.588200000000143
# Row 132
$ws.Range("H132").Value = 865.0714
$ws.Range("I132").Value = 753.7451
$ws.Range("J132").Value = 2000.6
$ws.Range("K132").Value = 2261.2353
$ws.Range("L132").Value = 6001.799999999999
$ws.Range("M132").Value = 268.7647000000002
$ws.Range("N132").Value = -11061.8
# Row 138
$ws.Range("H138").Value = 1777.33
$ws.Range("J138").Value = 1983.1143
$ws.Range("L138").Value = 5949.3429
$ws.Range("N138").Value = -16229.3429
# Row 141
$ws.Range("H141").Value = 1122273.9
$ws.Range("J141").Value = 3072.2727
$ws.Range("L141").Value = 9216.8181
$ws.Range("N141").Value = -19576.8181

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 35688.043
$ws.Range("J61").Value = 8787
$ws.Range("L61").Value = 8787
$ws.Range("N61").Value = -9211
# Row 74
$ws.Range("H74").Value = 735.3714
$ws.Range("I74").Value = 742.82355
$ws.Range("J74").Value = 482
$ws.Range("K74").Value = 742.82355
$ws.Range("L74").Value = 482
$ws.Range("M74").Value = 131.17645
$ws.Range("N74").Value = -2230
# Row 77
$ws.Range("H77").Value = 735.3714
$ws.Range("I77").Value = 742.82355
$ws.Range("J77").Value = 482
$ws.Range("K77").Value = 3714.11775
$ws.Range("L77").Value = 2410
$ws.Range("M77").Value = 653.8822500000001
$ws.Range("N77").Value = -11146
# Row 122
$ws.Range("H122").Value = 1334.1482
$ws.Range("I122").Value = 1305.0834
$ws.Range("K122").Value = 3915.2502
$ws.Range("M122").Value = -1465.2502
# Row 123
$ws.Range("H123").Value = 81499.5
$ws.Range("J123").Value = 81499.5
$ws.Range("L123").Value = 81499.5
$ws.Range("N123").Value = -91299.5
# Row 132
$ws.Range("H132").Value = 1688.6615
$ws.Range("I132").Value = 1173
$ws.Range("K132").Value = 3519
$ws.Range("M132").Value = -989
# Row 136
$ws.Range("H136").Value = 35688.043
$ws.Range("J136").Value = 8787
$ws.Range("L136").Value = 26361
$ws.Range("N136").Value = -31461

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3520.25
$ws.Range("I31").Value = 2685
$ws.Range("J31").Value = 4355.5
$ws.Range("K31").Value = 2685
$ws.Range("L31").Value = 4355.5
$ws.Range("M31").Value = -2390
$ws.Range("N31").Value = -4945.5
# Row 34
$ws.Range("H34").Value = 3520.25
$ws.Range("I34").Value = 2685
$ws.Range("J34").Value = 4355.5
$ws.Range("K34").Value = 2685
$ws.Range("L34").Value = 4355.5
$ws.Range("M34").Value = -2483
$ws.Range("N34").Value = -4759.5
# Row 134
$ws.Range("H134").Value = 1351.9117
$ws.Range("I134").Value = 1221.1852
$ws.Range("K134").Value = 3663.5556
$ws.Range("M134").Value = -1128.5556
# Row 141
$ws.Range("H141").Value = 58888.777
$ws.Range("J141").Value = 57249.875
$ws.Range("L141").Value = 57249.875
$ws.Range("N141").Value = -67609.875

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 604.53845
$ws.Range("J5").Value = 899.3333
$ws.Range("L5").Value = 2697.9999
$ws.Range("N5").Value = -2921.9999
# Row 36
$ws.Range("H36").Value = 1666.5
$ws.Range("J36").Value = 1399.8
$ws.Range("L36").Value = 4199.4
$ws.Range("N36").Value = -4537.4
# Row 131
$ws.Range("H131").Value = 16338.462
$ws.Range("J131").Value = 17664.688
$ws.Range("L131").Value = 52994.064
$ws.Range("N131").Value = -63074.064
# Row 135
$ws.Range("H135").Value = 604.53845
$ws.Range("J135").Value = 899.3333
$ws.Range("L135").Value = 8093.9997
$ws.Range("N135").Value = -13163.9997

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 1607.4166
$ws.Range("I102").Value = 1680.8182
$ws.Range("K102").Value = 1680.8182
$ws.Range("M102").Value = -58.81819999999993
# Row 110
$ws.Range("H110").Value = 66490
$ws.Range("J110").Value = 66490
$ws.Range("L110").Value = 66490
$ws.Range("N110").Value = -74670
# Row 122
$ws.Range("H122").Value = 1939.6
$ws.Range("I122").Value = 1233
$ws.Range("K122").Value = 3699
$ws.Range("M122").Value = -1249
# Row 132
$ws.Range("H132").Value = 840239.5600000001
$ws.Range("J132").Value = 4365.5713
$ws.Range("L132").Value = 13096.7139
$ws.Range("N132").Value = -18156.7139

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3261
$ws.Range("I7").Value = 3161.6
$ws.Range("K7").Value = 3161.6
$ws.Range("M7").Value = -3049.6
# Row 40
$ws.Range("H40").Value = 9195.84
$ws.Range("I40").Value = 8608.954
$ws.Range("K40").Value = 8608.954
$ws.Range("M40").Value = -8472.954
# Row 46
$ws.Range("H46").Value = 2369
$ws.Range("J46").Value = 2799.625
$ws.Range("L46").Value = 2799.625
$ws.Range("N46").Value = -3175.625
# Row 59
$ws.Range("H59").Value = 398
$ws.Range("J59").Value = 398
$ws.Range("L59").Value = 398
$ws.Range("N59").Value = -1706
# Row 93
$ws.Range("H93").Value = 1579.7894
$ws.Range("I93").Value = 953.5714
$ws.Range("J93").Value = 3333.2
$ws.Range("K93").Value = 953.5714
$ws.Range("L93").Value = 3333.2
$ws.Range("M93").Value = 294.4286
$ws.Range("N93").Value = -5829.2
# Row 109
$ws.Range("H109").Value = 49999
$ws.Range("J109").Value = 49999
$ws.Range("L109").Value = 49999
$ws.Range("N109").Value = -52773
# Row 126
$ws.Range("H126").Value = 3261
$ws.Range("I126").Value = 3161.6
$ws.Range("K126").Value = 9484.799999999999
$ws.Range("M126").Value = -7014.799999999999
# Row 136
$ws.Range("H136").Value = 2597.5625
$ws.Range("I136").Value = 2550.8462
$ws.Range("K136").Value = 7652.5386
$ws.Range("M136").Value = -5102.5386

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 94
$ws.Range("J94").Value = 29500
$ws.Range("L94").Value = 29500
$ws.Range("N94").Value = -31302
# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
# Row 132
$ws.Range("H132").Value = 1227.826
$ws.Range("I132").Value = 1120.7344
$ws.Range("K132").Value = 3362.2032
$ws.Range("M132").Value = -832.2031999999999
# Row 136
$ws.Range("H136").Value = 1804.2693
$ws.Range("I136").Value = 1308.8823
$ws.Range("K136").Value = 3926.6469
$ws.Range("M136").Value = -1376.6469
